$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 5: Quantity 2 -> 1, Total Cost 18.50 -> 9.25
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "1"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.25"

# Add new row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "004061"
$ws.Range("B10").Value = "Natalie's - Honey Tangerine"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "14.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.00"
